# "Generate Report for Handoff"
#
# The localization-status report previously had its two tracked files
# ("18633056-ff34-44e2-8461-f8cb6b01ebaa.md" and
#  "6100965e-3277-4e74-8ceb-b89abe4613f0.md") mixed up: row 2 carried the
# data that actually belonged to "6100965e..." and row 3 carried the data
# that actually belonged to "18633056...". This handoff-report generation
# pass fixes the row/file association and records a fresh handoff event
# (new status + timestamps) for "18633056...".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 now belongs to file 6100965e... (previously handed back, unchanged)
$ov.Range("A2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-03-19 16:45:13"

# Row 3 now belongs to file 18633056..., freshly handed off
$ov.Range("A3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-19 16:46:54"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-19 16:45:03"
$zh.Range("F2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.md"
$zh.Range("G2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-19 16:45:56"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-19 16:46:46"
$zh.Range("F3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$zh.Range("G3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-19 16:45:56"
$zh.Range("J3").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.de-de.xlf"
$de.Range("E2").Value = "2016-03-19 16:45:13"
$de.Range("F2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.md"
$de.Range("G2").Value = "6100965e-3277-4e74-8ceb-b89abe4613f0.3b6a021dba2f26b9dcc1733c6e892d2eeeb0f51a.de-de.xlf"
$de.Range("H2").Value = "2016-03-19 16:46:11"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf"
$de.Range("E3").Value = "2016-03-19 16:46:54"
$de.Range("F3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.md"
$de.Range("G3").Value = "18633056-ff34-44e2-8461-f8cb6b01ebaa.27ff3c630544d83dd9fe630cd8fb891bc74d74fa.de-de.xlf"
$de.Range("H3").Value = "2016-03-19 16:46:11"
$de.Range("J3").Value = "Include"
